$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# --- Text fixes in shared strings (content stays the same, just relabeled) ---
# Order matters for shared-string table indexing: new strings must be
# registered in the same order the reference workbook used.
$ws.Range("E4").Value = "(Xi-średnia)^2*ni"
$ws.Range("B3").Value = "Xi"

# --- Formula fixes: use odch.stand (C18) instead of wariancja (C17) ---
$ws.Range("C20").Formula = "=C16-C19*(C18/SQRT(C15))"
$ws.Range("C21").Formula = "=C16+C19*(C18/SQRT(C15))"

# --- Column L width ---
$ws.Range("L1").ColumnWidth = 42.42578125

# --- Selection change ---
$ws.Range("I13").Select()
